$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (High Pressure Alarm): action text now reflects a restart, not just continued ventilation ---
$ws.Range("C2").Value = "Audible/visual alarm; stop ventilation cycle and reset to try again"

# --- Row 6 (Disconnect Alarm): "when does it occur" now describes a low plateau pressure condition ---
$ws.Range("B6").Value = "Plateau pressure is below the set low plateau pressure limit"
# Touch the whole-cell font so a dedicated (plain) cell style is registered,
# then color the measured/set-point phrases per the red/purple convention used elsewhere.
$ws.Range("B6").Font.Name = "Calibri"
$ws.Range("B6").Font.Size = 11
$ws.Range("B6").Characters(1, 16).Font.Color = 255
$ws.Range("B6").Characters(34, 27).Font.Color = 10498160

# --- Column C is now wider to fit the longer action text ---
$ws.Columns("C").ColumnWidth = 62.3

# --- View state: zoomed to 85% with the selection left on C12 ---
$excel.ActiveWindow.Zoom = 85
$ws.Range("C12").Select()
